$d = $word.ActiveDocument
$t = $d.Tables(1)

# Each entry: row (1-based), column (1-based), new cell text.
# Row/column indices correspond to the physical <w:tc> cells in the
# single table that spans the whole document.
$changes = @(
    @{ Row = 3;  Col = 2; Text = "YesTemperature gradients are the driving force for heat transport. Thermal conductivity and heat capacity are temperature dependent." },
    @{ Row = 7;  Col = 5; Text = "NoBut indirectly through temperature" },
    @{ Row = 11; Col = 5; Text = "NoBut indirectly through temperature" },
    @{ Row = 15; Col = 5; Text = "NoBut indirectly through temperature" },
    @{ Row = 19; Col = 2; Text = "YesAffects heat flux from repository. Canister spacing particularly important in the near field." },
    @{ Row = 23; Col = 5; Text = "NoBut indirectly through rock stresses and temperature." },
    @{ Row = 27; Col = 5; Text = "NoBut indirectly through temperature." },
    @{ Row = 31; Col = 2; Text = "YesDetermines thermal properties." },
    @{ Row = 35; Col = 2; Text = "YesMarginally and locally." },
    @{ Row = 35; Col = 5; Text = "NoBut indirectly through temperature and groundwater composition." },
    @{ Row = 39; Col = 5; Text = "NoBut indirectly through temperature." },
    @{ Row = 51; Col = 2; Text = "YesAffects scope and extent of convective heat transport." },
    @{ Row = 51; Col = 5; Text = "NoBut, indirectly through temperature." }
)

foreach ($chg in $changes) {
    $cell = $t.Cell($chg.Row, $chg.Col)
    $rng = $cell.Range
    # Drop the trailing cell-mark character so we only overwrite the
    # visible text, keeping the cell's own formatting/run intact.
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $chg.Text
}
